$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark attendance ("1") for the additional day column (G) on rows 4, 6 and 7,
# matching the existing pattern already used for rows 2, 3 and 5.
$ws.Range("G4").Value = "1"
$ws.Range("G6").Value = "1"
$ws.Range("G7").Value = "1"

# Move the active selection to G7, as recorded after the edit.
$ws.Range("G7").Select()
